$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet updates
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("About")

# Insert 4 new rows right after row 3 for the expanded source citation
$ws.Range("A4:A7").EntireRow.Insert()

# Row 7: cited table reference (set first so the shared-string table ends up
# in the same order the source workbook uses)
$ws.Range("A7").Font.Bold = $true
$ws.Range("B7").Value = "Table 5 Generalized Cost Coefficient Calibration"

# Update the "Sources:" line - replace old note with new EPA citation header
$ws.Range("B3").Value = "United States EPA"

# Row 4: publication year (left aligned number)
$ws.Range("A4").Font.Bold = $true
$ws.Range("B4").Value = 2012
$ws.Range("B4").HorizontalAlignment = -4131

# Row 5: report title
$ws.Range("A5").Font.Bold = $true
$ws.Range("B5").Value = "Consumer Vehicle Choice Model Documentation"

# Row 6: report URL
$ws.Range("A6").Font.Bold = $true
$ws.Range("B6").Value = "https://nepis.epa.gov/Exe/ZyPDF.cgi/P100EZ37.PDF?Dockey=P100EZ37.PDF"

# Insert 3 new rows after the "data on technology buyers' behavior." line
# (old row 11, now shifted to row 15) to hold the new calibration explanation
$ws.Range("A16:A18").EntireRow.Insert()
$ws.Range("A16").Value = "We choose a value of -3 for passenger vehicles and a value of -5 for other vehicle types, "
$ws.Range("A17").Value = "based on the ranges in Table 5 of the cited EPA documentation."

# ---------------------------------------------------------------------------
# "TTLE" sheet updates - change logit exponent values from -3 to -5
# ---------------------------------------------------------------------------
$ts = $wb.Worksheets.Item("TTLE")
$ts.Range("B2:C7").Value = -5
